$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r, B (Coin), C (Link), D (Price), E (Volume 1h)
$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "28.402.24", "  +0.76%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.865.29", "  -0.75%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.017", "  +1.36%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "315.78", "  +0.85%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.016", "  +1.39%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.5120", "  -0.14%  "),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3918", "  +0.52%  "),
    @(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.08305", "  -0.71%  "),
    @(10, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.115", "  -0.19%  "),
    @(11, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "6.236", "  +0.08%  "),
    @(12, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.860.73", "  -0.78%  "),
    @(13, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "20.28", "  -2.14%  "),
    @(14, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.208", "  -1.25%  "),
    @(15, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.016", "  +1.37%  "),
    @(16, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.00001101", "  -0.50%  "),
    @(17, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "91.19", "  -0.17%  "),
    @(18, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06717", "  +1.02%  "),
    @(19, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "17.63", "  -0.67%  "),
    @(20, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.015", "  +1.32%  "),
    @(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.947", "  -1.70%  "),
    @(22, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "28.431.69", "  +0.76%  "),
    @(23, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.09", "  -1.06%  "),
    @(24, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.258", "  -0.28%  "),
    @(25, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.064.08", "  -1.29%  "),
    @(26, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "161.05", "  +1.60%  "),
    @(27, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "20.55", "  -0.39%  "),
    @(28, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.405", "  -4.26%  "),
    @(29, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "126.73", "  +1.07%  "),
    @(30, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1051", "  -1.28%  "),
    @(31, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.034", "  -0.81%  "),
    @(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.801", "  -1.47%  "),
    @(33, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "3.638", "  +1.39%  "),
    @(34, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02440", "  -0.65%  "),
    @(35, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "9.198", "  -5.41%  "),
    @(36, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.06485", "  -1.05%  "),
    @(37, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.2170", "  -1.12%  "),
    @(38, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.251", "  +1.53%  "),
    @(39, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.184", "  -2.22%  "),
    @(40, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.6421", "  -1.44%  "),
    @(41, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.961", "  -1.39%  "),
    @(42, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "11.11", "  -1.70%  "),
    @(43, "Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.6000", "  -1.87%  "),
    @(44, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "12.96", "  -1.04%  "),
    @(45, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.707", "  +0.90%  "),
    @(46, "WEMIXTOKEN", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.284", "  -0.56%  "),
    @(47, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.984", "  -1.64%  "),
    @(48, "EOS", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos", "1.202", "  -2.23%  "),
    @(49, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "120.95", "  -0.57%  "),
    @(50, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.06863", "  -0.57%  "),
    @(51, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "76.12", "  -2.70%  "),
)

# Rows whose Price (column D) text would otherwise be auto-parsed by Excel as a
# number; force these cells to Text format first so the literal string is kept.
$textPriceRows = @(4, 5, 6, 7, 8, 9, 10, 11, 13, 14, 15, 16, 17, 18, 19, 20, 21, 23, 24, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51)
foreach ($r in $textPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
